$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update remaining values: A2 -> 30, A3 -> 31
$ws.Range("A2").Value = 30
$ws.Range("A3").Value = 31

# Clear rows 4 through 6 (values 27, 28, 29) without shifting other cells,
# so the used range shrinks to A1:A3
$ws.Range("A4:A6").ClearContents()
